$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "2023 October", 64, 9),
    @(3, "2023 November", 64, 4),
    @(4, "2023 December", 55, 7),
    @(5, "2024 January", 53, 8),
    @(6, "2024 February", 53, 14),
    @(7, "2024 March", 34, 9),
    @(8, "2024 April", 34, 13),
    @(9, "2024 May", 54, 24),
    @(10, "2024 June", 37, 25),
    @(11, "2024 July", 46, 19),
    @(12, "2024 August", 34, 33),
    @(13, "2024 September", 23, 2)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $month = $entry[1]
    $closed = $entry[2]
    $opened = $entry[3]

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $month

    $ws.Cells.Item($row, 2).Value = $closed
    $ws.Cells.Item($row, 3).Value = $opened
}
